$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Build the repeated "TC_ID / TC_Description / Steps / Status" test-case
# block: a bold+bordered header row followed by a bordered data row, for
# test cases 1..8 (rows 4..19), plus one trailing header row (row 20).
# ---------------------------------------------------------------------------

$headerValues = @("TC_ID", "TC_Description", "Steps", "Status (Pass/Fail)")
$stepsText = "Navigate to expedia.com`nClick on check-in input field`nEnter date`n"

$row = 4
for ($tc = 1; $tc -le 8; $tc++) {

    # Header row: bold font + thin box border on A:D
    $hdr = $ws.Range("A$row`:D$row")
    for ($c = 1; $c -le 4; $c++) {
        $ws.Cells.Item($row, $c).Value = $headerValues[$c - 1]
    }
    $hdr.Font.Bold = $true
    $hdr.Borders.LineStyle = 1

    $dataRow = $row + 1

    # A: test case number, formatted as integer "0", with border
    $aCell = $ws.Cells.Item($dataRow, 1)
    $aCell.Value = $tc
    $aCell.NumberFormat = "0"
    $aCell.Borders.LineStyle = 1

    # B: short description, with border
    $bCell = $ws.Cells.Item($dataRow, 2)
    $bCell.Value = "Enter checkin date"
    $bCell.Borders.LineStyle = 1

    # C: multi-line steps, wrapped text, with border
    $cCell = $ws.Cells.Item($dataRow, 3)
    $cCell.Value = $stepsText
    $cCell.WrapText = $true
    $cCell.Borders.LineStyle = 1

    # D: empty status cell, with border only
    $dCell = $ws.Cells.Item($dataRow, 4)
    $dCell.Borders.LineStyle = 1

    $ws.Range("A$dataRow`:D$dataRow").RowHeight = 75

    $row = $row + 2
}

# Trailing header row (row 20) with no following data row
$hdr = $ws.Range("A$row`:D$row")
for ($c = 1; $c -le 4; $c++) {
    $ws.Cells.Item($row, $c).Value = $headerValues[$c - 1]
}
$hdr.Font.Bold = $true
$hdr.Borders.LineStyle = 1

# ---------------------------------------------------------------------------
# Column widths for the new B/C/D columns (closest values reachable through
# the ColumnWidth property, which snaps to the host's pixel grid).
# ---------------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 16.6
$ws.Columns.Item(3).ColumnWidth = 22.6
$ws.Columns.Item(4).ColumnWidth = 15.5

# Match the saved selection
$ws.Range("H6").Select() | Out-Null
